$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): strip the bold/border/center-top style that was
# applied to A1:T1, and clear the stray "Unnamed: 0" label in A1 ---
$ws.Range("A1:T1").Style = "Normal"
$ws.Range("A1").Value = ""

# --- Recalculated metrics for rows 3-7 (pre/post/total fixation data
# cleaning) and clearing the "param" (O) column which no longer has data ---

# Row 3: Revisit count
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 16
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 14
$ws.Range("J3").Value = 21
$ws.Range("K3").Value = 9
$ws.Range("L3").Value = 14
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = 25
$ws.Range("Q3").Value = 5
$ws.Range("S3").Value = 13
$ws.Range("T3").Value = 11

# Row 4: Fixation count
$ws.Range("B4").Value = 25
$ws.Range("C4").Value = 93
$ws.Range("D4").Value = 29
$ws.Range("E4").Value = 41
$ws.Range("F4").Value = 27
$ws.Range("J4").Value = 53
$ws.Range("K4").Value = 15
$ws.Range("L4").Value = 25
$ws.Range("M4").Value = 4
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = 124
$ws.Range("Q4").Value = 8
$ws.Range("S4").Value = 20
$ws.Range("T4").Value = 16

# Row 5: Dwell time (ms)
$ws.Range("B5").Value = 9426.370000000001
$ws.Range("C5").Value = 28496.49
$ws.Range("D5").Value = 10828.1
$ws.Range("E5").Value = 14048.1
$ws.Range("F5").Value = 10511.04
$ws.Range("J5").Value = 17634.8
$ws.Range("K5").Value = 5455.02
$ws.Range("L5").Value = 9426.370000000001
$ws.Range("M5").Value = 2836.2
$ws.Range("N5").Value = 1334.75
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = 32337.92
$ws.Range("Q5").Value = 3035.57
$ws.Range("S5").Value = 7608.08
$ws.Range("T5").Value = 5722.42

# Row 6: Dwell time (%)
$ws.Range("B6").Value = 6.32
$ws.Range("C6").Value = 19.09
$ws.Range("D6").Value = 7.26
$ws.Range("E6").Value = 9.41
$ws.Range("F6").Value = 7.04
$ws.Range("H6").Value = 0.1
$ws.Range("J6").Value = 11.82
$ws.Range("K6").Value = 3.66
$ws.Range("L6").Value = 6.32
$ws.Range("M6").Value = 1.9
$ws.Range("N6").Value = 0.89
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = 21.67
$ws.Range("Q6").Value = 2.03
$ws.Range("S6").Value = 5.1
$ws.Range("T6").Value = 3.83

# Row 7: Fixation duration (ms)
$ws.Range("B7").Value = 377.05
$ws.Range("C7").Value = 306.41
$ws.Range("D7").Value = 373.38
$ws.Range("E7").Value = 342.64
$ws.Range("F7").Value = 389.3
$ws.Range("J7").Value = 332.73
$ws.Range("K7").Value = 363.67
$ws.Range("L7").Value = 377.05
$ws.Range("M7").Value = 709.05
$ws.Range("N7").Value = 266.95
$ws.Range("O7").Value = ""
$ws.Range("P7").Value = 260.79
$ws.Range("Q7").Value = 379.45
$ws.Range("S7").Value = 380.4
$ws.Range("T7").Value = 357.65

# Row 8: First fixation duration (ms) - values unchanged, only the "param"
# column had data and it was cleared
$ws.Range("O8").Value = ""

# --- The last two blank padding rows (10-11) are no longer part of the
# cleaned data range, so remove them entirely and shift the dimension up ---
$ws.Range("A10:A11").EntireRow.Delete()
